# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) values were recomputed by the upstream data
# pipeline (std/mean regen + s_vals calc) and rewritten here with the
# freshly computed integers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 3
    4  = 1
    5  = 0
    6  = 0
    7  = 1
    8  = 1
    9  = 1
    10 = 2
    11 = 2
    12 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 0
    23 = 1
    24 = 0
    25 = 1
    26 = 0
    27 = 1
    28 = 1
    29 = 0
    30 = 1
    31 = 1
    32 = 2
    33 = 0
    34 = 2
    35 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
